# Refresh cached Universalis market-price snapshots (and the leve-profit
# figures derived from them) on the Typhon Profits workbook sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Cells.Item(17, 8).Value = 541.6111
$ws.Cells.Item(17, 10).Value = 538.1429000000001
$ws.Cells.Item(17, 12).Value = 1614.4287
$ws.Cells.Item(17, 14).Value = -1950.4287
# Row 38: Just Give Him a Serum
$ws.Cells.Item(38, 8).Value = 161.8
$ws.Cells.Item(38, 9).Value = 161.8
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 485.4
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = -113.4
$ws.Cells.Item(38, 14).ClearContents()
# Row 53: No Accounting for Waste
$ws.Cells.Item(53, 8).Value = 3748.889
$ws.Cells.Item(53, 9).Value = 383
$ws.Cells.Item(53, 11).Value = 383
$ws.Cells.Item(53, 13).Value = 254
# Row 70: Consecrating Congregation
$ws.Cells.Item(70, 8).Value = 1407.2
$ws.Cells.Item(70, 9).Value = 1196
$ws.Cells.Item(70, 10).Value = 1900
$ws.Cells.Item(70, 11).Value = 3588
$ws.Cells.Item(70, 12).Value = 5700
$ws.Cells.Item(70, 13).Value = -3318
$ws.Cells.Item(70, 14).Value = -6240
# Row 73: Curbing the Contagion (L)
$ws.Cells.Item(73, 8).Value = 1407.2
$ws.Cells.Item(73, 9).Value = 1196
$ws.Cells.Item(73, 10).Value = 1900
$ws.Cells.Item(73, 11).Value = 3588
$ws.Cells.Item(73, 12).Value = 5700
$ws.Cells.Item(73, 13).Value = -2652
$ws.Cells.Item(73, 14).Value = -7572
# Row 86: Filling in the Blanks
$ws.Cells.Item(86, 8).Value = 9274.571
$ws.Cells.Item(86, 9).Value = 2314.5715
$ws.Cells.Item(86, 10).Value = 16234.571
$ws.Cells.Item(86, 11).Value = 2314.5715
$ws.Cells.Item(86, 12).Value = 16234.571
$ws.Cells.Item(86, 13).Value = -1191.5715
$ws.Cells.Item(86, 14).Value = -18480.571
# Row 89: Ink into Antiquity (L)
$ws.Cells.Item(89, 8).Value = 9274.571
$ws.Cells.Item(89, 9).Value = 2314.5715
$ws.Cells.Item(89, 10).Value = 16234.571
$ws.Cells.Item(89, 11).Value = 11572.8575
$ws.Cells.Item(89, 12).Value = 81172.855
$ws.Cells.Item(89, 13).Value = -5956.8575
$ws.Cells.Item(89, 14).Value = -92404.855
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 3120.7144
$ws.Cells.Item(32, 9).Value = 2091.946
$ws.Cells.Item(32, 11).Value = 2091.946
$ws.Cells.Item(32, 13).Value = -1804.946
$ws = $wb.Worksheets.Item("CRP")
# Row 41: The Lone Bowman
$ws.Cells.Item(41, 8).Value = 26440
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 26440
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 26440
$ws.Cells.Item(41, 13).ClearContents()
$ws.Cells.Item(41, 14).Value = -27296
# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 25031.191
$ws.Cells.Item(58, 9).Value = 1255.5264
$ws.Cells.Item(58, 11).Value = 1255.5264
$ws.Cells.Item(58, 13).Value = -1052.5264
# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 24069.5
$ws.Cells.Item(132, 9).Value = 36902.715
$ws.Cells.Item(132, 10).Value = 6103
$ws.Cells.Item(132, 11).Value = 110708.145
$ws.Cells.Item(132, 12).Value = 18309
$ws.Cells.Item(132, 13).Value = -108178.145
$ws.Cells.Item(132, 14).Value = -23369
# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 25031.191
$ws.Cells.Item(136, 9).Value = 1255.5264
$ws.Cells.Item(136, 11).Value = 3766.5792
$ws.Cells.Item(136, 13).Value = -1216.5792
$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face
$ws.Cells.Item(68, 8).Value = 1347.0588
$ws.Cells.Item(68, 10).Value = 1443.3334
$ws.Cells.Item(68, 12).Value = 4330.0002
$ws.Cells.Item(68, 14).Value = -5952.0002
# Row 70: Persona non Gratin
$ws.Cells.Item(70, 8).Value = 3298.5386
$ws.Cells.Item(70, 9).Value = 2226.375
$ws.Cells.Item(70, 11).Value = 6679.125
$ws.Cells.Item(70, 13).Value = -6364.125
# Row 71: No Margarine of Error (L)
$ws.Cells.Item(71, 8).Value = 1347.0588
$ws.Cells.Item(71, 10).Value = 1443.3334
$ws.Cells.Item(71, 12).Value = 12990.0006
$ws.Cells.Item(71, 14).Value = -21102.0006
# Row 73: Recipe for Disaster (L)
$ws.Cells.Item(73, 8).Value = 3298.5386
$ws.Cells.Item(73, 9).Value = 2226.375
$ws.Cells.Item(73, 11).Value = 6679.125
$ws.Cells.Item(73, 13).Value = -5587.125
# Row 76: Old Victories, New Tastes
$ws.Cells.Item(76, 8).Value = 4330.4546
$ws.Cells.Item(76, 9).Value = 1250
$ws.Cells.Item(76, 11).Value = 3750
$ws.Cells.Item(76, 13).Value = -3367
# Row 79: The Eats of Authenticity (L)
$ws.Cells.Item(79, 8).Value = 4330.4546
$ws.Cells.Item(79, 9).Value = 1250
$ws.Cells.Item(79, 11).Value = 3750
$ws.Cells.Item(79, 13).Value = -2424
# Row 80: Saucy for a Suitor
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).ClearContents()
# Row 81: It Goes Down Smoothly
$ws.Cells.Item(81, 8).Value = 4025.7273
$ws.Cells.Item(81, 9).Value = 606.5
$ws.Cells.Item(81, 10).Value = 4785.5557
$ws.Cells.Item(81, 11).Value = 1819.5
$ws.Cells.Item(81, 12).Value = 14356.6671
$ws.Cells.Item(81, 13).Value = -696.5
$ws.Cells.Item(81, 14).Value = -16602.6671
# Row 82: Persuasion of a Higher Power
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()
# Row 83: Saved by the Sauce (L)
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).ClearContents()
# Row 84: Quenching the Flame (L)
$ws.Cells.Item(84, 8).Value = 4025.7273
$ws.Cells.Item(84, 9).Value = 606.5
$ws.Cells.Item(84, 10).Value = 4785.5557
$ws.Cells.Item(84, 11).Value = 5458.5
$ws.Cells.Item(84, 12).Value = 43070.0013
$ws.Cells.Item(84, 13).Value = 157.5
$ws.Cells.Item(84, 14).Value = -54302.0013
# Row 85: Loaves and Fishes (L)
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()
# Row 86: Let's Not Get Sappy
$ws.Cells.Item(86, 8).Value = 858.6667
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).ClearContents()
# Row 87: Soup That Eats Like a Knight
$ws.Cells.Item(87, 8).Value = 5280
$ws.Cells.Item(87, 9).Value = 636
$ws.Cells.Item(87, 10).Value = 28500
$ws.Cells.Item(87, 11).Value = 1908
$ws.Cells.Item(87, 12).Value = 85500
$ws.Cells.Item(87, 13).Value = -660
$ws.Cells.Item(87, 14).Value = -87996
# Row 88: Don't Let It Fall Apart
$ws.Cells.Item(88, 8).Value = 7743.2
$ws.Cells.Item(88, 10).Value = 7743.2
$ws.Cells.Item(88, 12).Value = 23229.6
$ws.Cells.Item(88, 14).Value = -24085.6
# Row 89: Luxury Spillover (L)
$ws.Cells.Item(89, 8).Value = 858.6667
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).ClearContents()
# Row 90: Like Ma Used to Make (L)
$ws.Cells.Item(90, 8).Value = 5280
$ws.Cells.Item(90, 9).Value = 636
$ws.Cells.Item(90, 10).Value = 28500
$ws.Cells.Item(90, 11).Value = 5724
$ws.Cells.Item(90, 12).Value = 256500
$ws.Cells.Item(90, 13).Value = 516
$ws.Cells.Item(90, 14).Value = -268980
# Row 91: Better Come Back with a Sandwich (L)
$ws.Cells.Item(91, 8).Value = 7743.2
$ws.Cells.Item(91, 10).Value = 7743.2
$ws.Cells.Item(91, 12).Value = 23229.6
$ws.Cells.Item(91, 14).Value = -26193.6
# Row 93: Loquacious
$ws.Cells.Item(93, 8).Value = 2239.3635
$ws.Cells.Item(93, 9).Value = 1816.5
$ws.Cells.Item(93, 10).Value = 2333.3333
$ws.Cells.Item(93, 11).Value = 5449.5
$ws.Cells.Item(93, 12).Value = 6999.999899999999
$ws.Cells.Item(93, 13).Value = -3577.5
$ws.Cells.Item(93, 14).Value = -10743.9999
# Row 131: The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 814.85
$ws.Cells.Item(131, 10).Value = 831.7292
$ws.Cells.Item(131, 12).Value = 2495.1876
$ws.Cells.Item(131, 14).Value = -12575.1876
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 3657.5715
$ws.Cells.Item(80, 9).Value = 2959.1667
$ws.Cells.Item(80, 10).Value = 4588.778
$ws.Cells.Item(80, 11).Value = 2959.1667
$ws.Cells.Item(80, 12).Value = 4588.778
$ws.Cells.Item(80, 13).Value = -1961.1667
$ws.Cells.Item(80, 14).Value = -6584.778
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 3657.5715
$ws.Cells.Item(83, 9).Value = 2959.1667
$ws.Cells.Item(83, 10).Value = 4588.778
$ws.Cells.Item(83, 11).Value = 14795.8335
$ws.Cells.Item(83, 12).Value = 22943.89
$ws.Cells.Item(83, 13).Value = -9803.833500000001
$ws.Cells.Item(83, 14).Value = -32927.89
# Row 113: Copious Crystal Cannons
$ws.Cells.Item(113, 8).Value = 2861.8572
$ws.Cells.Item(113, 9).Value = 2817.647
$ws.Cells.Item(113, 10).Value = 3049.75
$ws.Cells.Item(113, 11).Value = 2817.647
$ws.Cells.Item(113, 12).Value = 3049.75
$ws.Cells.Item(113, 13).Value = -647.6469999999999
$ws.Cells.Item(113, 14).Value = -7389.75
$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Cells.Item(96, 8).Value = 4812.5
$ws.Cells.Item(96, 10).Value = 5666.6665
$ws.Cells.Item(96, 12).Value = 5666.6665
$ws.Cells.Item(96, 14).Value = -8412.666499999999
# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 1889.862
$ws.Cells.Item(132, 9).Value = 1767.5294
$ws.Cells.Item(132, 10).Value = 2063.1667
$ws.Cells.Item(132, 11).Value = 5302.5882
$ws.Cells.Item(132, 12).Value = 6189.500100000001
$ws.Cells.Item(132, 13).Value = -2772.5882
$ws.Cells.Item(132, 14).Value = -11249.5001

Write-Host "Applied Typhon Profits market data refresh."
